$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9609120521172638
$ws.Range("C2").Value = 0.749185667752443

$ws.Range("B3").Value = 0.9609120521172638
$ws.Range("C3").Value = 0.758957654723127

$ws.Range("B4").Value = 0.9739413680781759
$ws.Range("C4").Value = 0.739413680781759

$ws.Range("B5").Value = 0.9674267100977199
$ws.Range("C5").Value = 0.758957654723127

$ws.Range("B6").Value = 0.9641693811074918
$ws.Range("C6").Value = 0.739413680781759
